# Insert a new "What about Production?" slide right before the closing
# "Q&A" slide (i.e. as the new 11th slide, pushing "Q&A" down to position 12).
#
# The deck currently has 11 slides, with "Q&A" as slide 11. Slides.Add(index,
# layout) inserts a brand new slide AT that position (shifting the old
# occupant of that position, and everything after it, one slot later) - so
# adding at index 11 lands the new slide directly in front of "Q&A".
#
# Layout 2 == ppLayoutText ("Title and Content"), the same layout used by the
# neighbouring slides: a Title placeholder plus a body/content placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Add(11, 2)

# --- Title placeholder -----------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
$title.Text = "What about Production?"
$title.Font.LanguageID = 1033

# --- Body / content placeholder --------------------------------------------
# One bullet per paragraph. Setting the language right after each paragraph
# is created (rather than once at the end, over the whole multi-paragraph
# range) is what gets it to actually stick on every paragraph's run.
$body = $s.Shapes.Item(2).TextFrame.TextRange
$lines = @("Scaling server side?", "JavaScript -> TypeScript", "Express -> Nest.js", "+ Unit tests", "+ React Component tests", "+ E2E test")

$body.Text = $lines[0]
$body.Paragraphs(1).Font.LanguageID = 1033

for ($i = 1; $i -lt $lines.Count; $i++) {
    $body.Text = $body.Text + "`r" + $lines[$i]
    $body.Paragraphs($i + 1).Font.LanguageID = 1033
}
